$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp header
$ws.Range("A1").Value = "Datos actualizados a 13 de Septiembre de 2020 a las 18:34"

# Update country case-count rows with refreshed data (values only; some
# rows swap position with a neighbour because the table is sorted by
# "Casos totales" descending and the refreshed totals changed the order)
# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 6690810
$ws.Range("C4").Value = 14209
$ws.Range("D4").Value = 3952366
$ws.Range("E4").Value = 2540205
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 111
$ws.Range("H4").Value = 198239

# Row 5: India
$ws.Range("A5").Value = "India"
$ws.Range("B5").Value = 4811712
$ws.Range("C5").Value = 59924
$ws.Range("D5").Value = 3749851
$ws.Range("E5").Value = 982637
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 610
$ws.Range("H5").Value = 79224

# Row 17: Reino Unido
$ws.Range("A17").Value = "Reino Unido"
$ws.Range("B17").Value = 368504
$ws.Range("C17").Value = 3330
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 41628

# Row 25: Alemania
$ws.Range("A25").Value = "Alemania"
$ws.Range("B25").Value = 260826
$ws.Range("C25").Value = 280
$ws.Range("D25").Value = 235300
$ws.Range("E25").Value = 16099
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 9427

# Row 34: Republica Dominicana
$ws.Range("A34").Value = "Republica Dominicana"
$ws.Range("B34").Value = 103660
$ws.Range("C34").Value = 568
$ws.Range("D34").Value = 77182
$ws.Range("E34").Value = 24510
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 15
$ws.Range("H34").Value = 1968

# Row 35: Rumania
$ws.Range("A35").Value = "Rumania"
$ws.Range("B35").Value = 103495
$ws.Range("C35").Value = 1109
$ws.Range("D35").Value = 43025
$ws.Range("E35").Value = 56307
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 36
$ws.Range("H35").Value = 4163

# Row 51: Etiopia
$ws.Range("A51").Value = "Etiopia"
$ws.Range("B51").Value = 64301
$ws.Range("C51").Value = 413
$ws.Range("D51").Value = 24983
$ws.Range("E51").Value = 38305
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 17
$ws.Range("H51").Value = 1013

# Row 52: Portugal
$ws.Range("A52").Value = "Portugal"
$ws.Range("B52").Value = 63983
$ws.Range("C52").Value = 673
$ws.Range("D52").Value = 44069
$ws.Range("E52").Value = 18047
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 7
$ws.Range("H52").Value = 1867

# Row 69: Chequia
$ws.Range("A69").Value = "Chequia"
$ws.Range("B69").Value = 35933
$ws.Range("C69").Value = 532
$ws.Range("D69").Value = 21209
$ws.Range("E69").Value = 14268
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 3
$ws.Range("H69").Value = 456

# Row 91: Grecia
$ws.Range("A91").Value = "Grecia"
$ws.Range("B91").Value = 13240
$ws.Range("C91").Value = 204
$ws.Range("D91").Value = 3804
$ws.Range("E91").Value = 9131
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 3
$ws.Range("H91").Value = 305

# Row 94: Albania
$ws.Range("A94").Value = "Albania"
$ws.Range("B94").Value = 11353
$ws.Range("C94").Value = 168
$ws.Range("D94").Value = 6569
$ws.Range("E94").Value = 4450
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 4
$ws.Range("H94").Value = 334

# Row 137: Trinidad yTobago
$ws.Range("A137").Value = "Trinidad yTobago"
$ws.Range("B137").Value = 3019
$ws.Range("C137").Value = 26
$ws.Range("D137").Value = 772
$ws.Range("E137").Value = 2194
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 2
$ws.Range("H137").Value = 53

# Row 138: Aruba
$ws.Range("A138").Value = "Aruba"
$ws.Range("B138").Value = 2994
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 1542
$ws.Range("E138").Value = 1434
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 18

# Row 140: Mali
$ws.Range("A140").Value = "Mali"
$ws.Range("B140").Value = 2924
$ws.Range("C140").Value = 8
$ws.Range("D140").Value = 2285
$ws.Range("E140").Value = 511
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 128

# Row 141: Reunion
$ws.Range("A141").Value = "Reunion"
$ws.Range("B141").Value = 2805
$ws.Range("C141").Value = 82
$ws.Range("D141").Value = 1313
$ws.Range("E141").Value = 1477
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 1
$ws.Range("H141").Value = 15

# Row 142: Birmania
$ws.Range("A142").Value = "Birmania"
$ws.Range("B142").Value = 2796
$ws.Range("C142").Value = 201
$ws.Range("D142").Value = 676
$ws.Range("E142").Value = 2104
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 16

# Row 151: Sierra Leona
$ws.Range("A151").Value = "Sierra Leona"
$ws.Range("B151").Value = 2109
$ws.Range("C151").Value = 13
$ws.Range("D151").Value = 1636
$ws.Range("E151").Value = 401
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 72

